$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "The cluster center can be characterised as: relatively high on outgoing, normal on grounded, normal on modest, conscientious, and considerate, normal on conceptual, normal on meticulous, self-reliant"
$ws.Range("A3").Value = "The cluster center can be characterised as: normal on outgoing, normal on pensive, relatively high on outgoing, candid, and conventional, slightly high on empathic, normal on casual"
$ws.Range("A4").Value = "The cluster center can be characterised as: slightly low on Reserved, very high on grounded, normal on modest, conscientious, and considerate, relatively high on empathic, normal on casual"
